## Preparing to add images for Numbers
## Adds an "Image Filename" value (0.png .. 10.png) in column I for the
## "Numbers" rows (zero..ten, rows 7-17), and moves the active selection
## to I18 (just past the newly-filled block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$imageNames = @(
    "0.png",
    "1.png",
    "2.png",
    "3.png",
    "4.png",
    "5.png",
    "6.png",
    "7.png",
    "8.png",
    "9.png",
    "10.png"
)

$startRow = 7
for ($i = 0; $i -lt $imageNames.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 9).Value = $imageNames[$i]
}

# Best-effort: line up the app window position too (not always reachable
# through this COM surface, so failures here are swallowed).
try {
    $win = $excel.ActiveWindow
    $win.Left = -120
    $win.Top = -120
} catch {
}

$ws.Activate()
$ws.Range("I18").Select()
